$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 791 (pushes existing rows 791..832 down to 793..834)
$ws.Range("A791:A792").EntireRow.Insert()

# Populate the two newly inserted rows with the new data point
# (2026/02/06, Friday, hour 18, rank 201) and (2026/02/06, Friday, hour 22, rank 201)
# Force the date column to Text format first so "2026/02/06" is stored as a
# literal string (matching the rest of the column) instead of being
# auto-parsed into a date serial number; then reset the style back to
# Normal so no stray number-format style sticks to the cell.
$ws.Range("A791").NumberFormat = "@"
$ws.Range("A791").Value = "2026/02/06"
$ws.Range("A791").Style = "Normal"
$ws.Range("B791").Value = "金"
$ws.Range("C791").Value = 18
$ws.Range("D791").Value = 201

$ws.Range("A792").NumberFormat = "@"
$ws.Range("A792").Value = "2026/02/06"
$ws.Range("A792").Style = "Normal"
$ws.Range("B792").Value = "金"
$ws.Range("C792").Value = 22
$ws.Range("D792").Value = 201
